$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# ManageNewsPage
# ---------------------------------------------------------------
$wsNews = $wb.Worksheets.Item("ManageNewsPage")
$wsNews.Range("A1").Value = "Original Text"
$wsNews.Range("B1").Value = "Updation Text"
$wsNews.Range("A2").Value = "Hello, Selenium Test"
$wsNews.Range("B2").Value = "Hello, TestNG Test"
$wsNews.Columns.Item(1).ColumnWidth = 23.7265625
$wsNews.Columns.Item(2).ColumnWidth = 17.54296875
$wsNews.Range("B2").Select()

# ---------------------------------------------------------------
# ManageFooterTextPage
# ---------------------------------------------------------------
$wsFooter = $wb.Worksheets.Item("ManageFooterTextPage")
$wsFooter.Range("A2").Value = "FNmbr: 1207B, Confident Atlenna"
$wsFooter.Range("A4").Value = 9876512345
$wsFooter.Columns.Item(1).ColumnWidth = 30.81640625
$wsFooter.Range("D5").Select()

# ---------------------------------------------------------------
# ManageCategory
# ---------------------------------------------------------------
$wsCategory = $wb.Worksheets.Item("ManageCategory")
$wsCategory.Range("A1").Value = "teststaff123"
$wsCategory.Range("A2").Value = "bunt"
$wsCategory.Range("A2").Select()

# ---------------------------------------------------------------
# AdminUsers
# ---------------------------------------------------------------
$wsAdmin = $wb.Worksheets.Item("AdminUsers")
$wsAdmin.Range("A5").Value = "testdb2"
$wsAdmin.Range("B5").Value = "testdb2"
$wsAdmin.Range("B5").Select()

# ---------------------------------------------------------------
# ManageContactPage
# ---------------------------------------------------------------
$wsContact = $wb.Worksheets.Item("ManageContactPage")
$wsContact.Range("A2").Value = 9945995541
$wsContact.Range("A4").Value = "Flat No 1234"
$wsContact.Range("A5").Value = 75
$wsContact.Range("A6").Value = 15
$wsContact.Range("A2").Select()

# Restore the originally active sheet/tab (ManageContactPage, activeTab=5)
$wsContact.Activate()
